$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Clear out the old "F:I" helper block (labels + regression parameters)
#    before we rebuild the sheet in its new column layout (J:O).
# ---------------------------------------------------------------------------
$ws.Range("F1:I21").UnMerge()
$ws.Range("F1:I21").Clear()

# ---------------------------------------------------------------------------
# 2. Header row (row 1)
# ---------------------------------------------------------------------------
$ws.Range("C1").Value = "Previsto - 10 primeiros"
$ws.Range("D1").Value = "Erro"
$ws.Range("E1").Value = "Previsto - 5 pri. 5 ult."
$ws.Range("F1").Value = "Erro"
$ws.Range("G1").Value = "Previsto - todos"
$ws.Range("H1").Value = "Erro"

# ---------------------------------------------------------------------------
# 3. Data rows 2-21: predictions + error columns
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 21; $r++) {
    $ws.Range("C$r").Formula = "=`$K`$2+`$M`$2*A$r"
    $ws.Range("D$r").Formula = "=B$r-C$r"
    $ws.Range("E$r").Formula = "=`$K`$4+`$M`$4*A$r"
    $ws.Range("F$r").Formula = "=B$r-E$r"
    $ws.Range("G$r").Formula = "=`$K`$6+`$M`$6*A$r"
    $ws.Range("H$r").Formula = "=C$r-G$r"
}

# ---------------------------------------------------------------------------
# 4. Regression-parameter helper blocks (J:O), rows 1-6
# ---------------------------------------------------------------------------
# Block 1: "10 primeiros" (rows 1-2)
$ws.Range("J1").Value = "10 primeiros"
$ws.Range("J2").Value = "b0 = "
$ws.Range("K2").Formula = "=AVERAGE(B2:B11)-M2*AVERAGE(A2:A11)"
$ws.Range("L2").Value = "b1 ="
$ws.Range("M2").Formula = "=SUMPRODUCT(A2:A11,B2:B11)/SUMSQ(A2:A11)"
$ws.Range("N1").Value = "ERRO QM"
$ws.Range("N2").Formula = "=SUMSQ(D2:D11)"
$ws.Range("O2").Formula = "=SUMSQ(D12:D21)"

# Block 2: "5 primeiros e 5 ultimos" (rows 3-4)
$ws.Range("J3").Value = "5 primeiros e 5 últimos"
$ws.Range("J4").Value = "b0 = "
$ws.Range("K4").Formula = "=AVERAGE(B2:B6,B16:B21)-M4*AVERAGE(A2:A6,A16:A21)"
$ws.Range("L4").Value = "b1 ="
$ws.Range("M4").Formula = "=(SUMPRODUCT(A2:A6,B2:B6)+SUMPRODUCT(A16:A21,B16:B21))/(SUMSQ(A2:A6)+SUMSQ(A16:A21))"
$ws.Range("N4").Formula = "=SUMSQ(F2:F6)+SUMSQ(F17:F21)"
$ws.Range("O4").Formula = "=SUMSQ(F7:F16)"

# Block 3: "todos" (rows 5-6)
$ws.Range("J5").Value = "todos"
$ws.Range("J6").Value = "b0 = "
$ws.Range("K6").Formula = "=AVERAGE(B4:B23)-M6*AVERAGE(A4:A23)"
$ws.Range("L6").Value = "b1 ="
$ws.Range("M6").Formula = "=SUMPRODUCT(A4:A23,B4:B23)/SUMSQ(A4:A23)"
$ws.Range("N6").Formula = "=SUMSQ(H2:H21)"

# ---------------------------------------------------------------------------
# 5. Row 22 - total squared-error summary
# ---------------------------------------------------------------------------
$ws.Range("B22").Value = "Erro quadrático médio total:"
$ws.Range("D22").Formula = "=SUMSQ(D2:D21)"
$ws.Range("F22").Formula = "=SUMSQ(F2:F21)"
$ws.Range("H22").Formula = "=SUMSQ(H2:H21)"

# ---------------------------------------------------------------------------
# 6. Merges
# ---------------------------------------------------------------------------
$ws.Range("J1:M1").Merge()
$ws.Range("J3:M3").Merge()
$ws.Range("J5:M5").Merge()
$ws.Range("B22:C22").Merge()

Write-Host "stage1 ok"
